$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25 (index 24 in the "N" column) was an empty template row; fill it in
# with a new client record, matching the formatting already used by the
# rows directly above it (rows 23-24).

$ws.Range("B25").Value = "DANIELE LOPES"
$ws.Range("C25").Value = "461c508d51e4a2193a8c7a9335877607"

$ws.Range("D25").NumberFormat = "yyyy-mm-dd"
$ws.Range("D25").Value = 44855

$ws.Range("E25").NumberFormat = "#,##0"
$ws.Range("E25").Value = 365

$ws.Range("F25").Value = "-"
$ws.Range("G25").Value = "VENDA 21 (20/10)"
